$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row cells for columns M, N, O
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Match the header formatting used by the existing header cells (B1:L1)
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

# Fill data rows 2 through 34 with new column values
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"      # column M
    $ws.Cells.Item($r, 14).Value = 20140882     # column N
    $ws.Cells.Item($r, 15).Value = 0            # column O
}
